$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-categorize a handful of existing "Citizen/stakeholder engagement"
#    rows to "Conference or active meeting" (per commit message: "change
#    some stakeholders to presentations").
# ---------------------------------------------------------------------------
$rowsToRecategorize = @(143, 144, 193, 215, 216, 217, 218)
foreach ($r in $rowsToRecategorize) {
    $ws.Cells.Item($r, 3).Value = "Conference or active meeting"
}

# ---------------------------------------------------------------------------
# 2) Append two new survey response rows to the bottom of Table1, growing
#    the table (and sheet data) from A1:G239 to A1:G241.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Row 240 -----------------------------------------------------------
$ws.Cells.Item(239, 1).Copy($ws.Cells.Item(240, 1))
$ws.Cells.Item(240, 1).Value = 43818.635567129626
$ws.Cells.Item(240, 2).Value = "ahappel@sheddaquarium.org"
$ws.Cells.Item(240, 3).Value = "Media opportunity"
$ws.Cells.Item(239, 4).Copy($ws.Cells.Item(240, 4))
$ws.Cells.Item(240, 4).Value = 43790
$ws.Cells.Item(240, 5).Value = "http://iaglr.org/ll/2019-3-Fall_LL3.pdf"
$ws.Cells.Item(240, 7).Value = "Non-peer reviewed publication I wrote about how the black spot disease paper came to be. Also speaks on Great Lakes Fish Finder and iNaturalist"

# --- Row 241 -----------------------------------------------------------
$ws.Cells.Item(239, 1).Copy($ws.Cells.Item(241, 1))
$ws.Cells.Item(241, 1).Value = 43818.636030092595
$ws.Cells.Item(241, 2).Value = "ahappel@sheddaquarium.org"
$ws.Cells.Item(241, 3).Value = "Field research"
$ws.Cells.Item(239, 4).Copy($ws.Cells.Item(241, 4))
$ws.Cells.Item(241, 4).Value = 43817
$ws.Cells.Item(241, 6).Value = 1
$ws.Cells.Item(241, 7).Value = "Haerther Work Day - buckthorn removal "

# ---------------------------------------------------------------------------
# 3) Update the saved view state: scroll/selection now rests on C172
#    instead of the prior bottom-of-sheet selection at E239.
# ---------------------------------------------------------------------------
$ws.Range("C172").Select() | Out-Null
